$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 4, 6, 7, 9, 10, 12, 13, 15, 16)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = "/api" + $cell.Value2
}
